$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/22/2024  Through  1/28/2024"

# --- Crime complaint table updates (rows 15-29) ---
$ws.Range("G14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("H14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 3
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -18.75
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = -18.75
$ws.Range("L16").Value = 62.5
$ws.Range("M16").Value = -38.095238095238
$ws.Range("N16").Value = -87.850467289719
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 16.666666666666
$ws.Range("I17").Value = 21
$ws.Range("J17").Value = 18
$ws.Range("K17").Value = 16.666666666666
$ws.Range("L17").Value = 23.529411764705
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = -47.5
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 44
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 83.333333333333
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 24
$ws.Range("K18").Value = 83.333333333333
$ws.Range("L18").Value = 29.411764705882
$ws.Range("M18").Value = 37.5
$ws.Range("N18").Value = -50.561797752809
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -43.478260869565
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 67
$ws.Range("H19").Value = -41.791044776119
$ws.Range("I19").Value = 39
$ws.Range("J19").Value = 67
$ws.Range("K19").Value = -41.791044776119
$ws.Range("L19").Value = 2.631578947368
$ws.Range("M19").Value = 18.181818181818
$ws.Range("N19").Value = 21.875
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 9
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -18.181818181818
$ws.Range("M20").Value = -25
$ws.Range("N20").Value = -84.745762711864
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = 2.564102564102
$ws.Range("F21").Value = 126
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = -8.695652173913
$ws.Range("I21").Value = 126
$ws.Range("J21").Value = 138
$ws.Range("K21").Value = -8.695652173913
$ws.Range("L21").Value = 14.545454545454
$ws.Range("M21").Value = 18.867924528301
$ws.Range("N21").Value = -62.275449101796
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 3
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = -25
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -7.692307692307
$ws.Range("I23").Value = 12
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = -7.692307692307
$ws.Range("L23").Value = 20
$ws.Range("M23").Value = 33.333333333333
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -32.258064516129
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = -3.370786516853
$ws.Range("I24").Value = 86
$ws.Range("J24").Value = 89
$ws.Range("K24").Value = -3.370786516853
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -7.52688172043
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 180
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = 37.5
$ws.Range("I25").Value = 44
$ws.Range("J25").Value = 32
$ws.Range("K25").Value = 37.5
$ws.Range("L25").Value = 37.5
$ws.Range("M25").Value = 33.333333333333
$ws.Range("G14").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 1
$ws.Range("G14").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("H14").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 0
$ws.Range("G14").Copy($ws.Range("F26"))
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("G14").Copy($ws.Range("I26"))
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = -75
$ws.Range("L26").Value = -50
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 3
$ws.Range("L27").Value = 0
$ws.Range("H14").Copy($ws.Range("L28"))
$ws.Range("L28").Value = -100
$ws.Range("H14").Copy($ws.Range("L29"))
$ws.Range("L29").Value = -100
